$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @(22178843, 22178844, 22178845, 22178846, 22178847, 22178848, 22178849, 22178850, 22178851)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newIds[$i]
}

$ws.Range("B2:B10").Select()
